# Append the new scrape batch (rows 76-84) collected 2025-05-09 to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  @{ r = 76; a = "A043"; b = "경기도_광명시";   c = "https://www.gm.go.kr/pt/user/nftcBbs/BD_selectNftcBbsList.do?q_nftcBbsCode=1001&q_rowPerPage=90"; d = "영회원 수변공원 조성공사 내 전망대 신기술・특허공법 선정위원회 결과 공개"; e = 45784; f = 45786.369710648149 },
  @{ r = 77; a = "A043"; b = "경기도_광명시";   c = "https://www.gm.go.kr/pt/user/nftcBbs/BD_selectNftcBbsList.do?q_nftcBbsCode=1001&q_rowPerPage=90"; d = "2025년 하반기 광명시 여성비전센터 단기특강 제안서 공모"; e = 45784; f = 45786.369710648149 },
  @{ r = 78; a = "A047"; b = "경기도_김포시";   c = "https://www.gimpo.go.kr/portal/ntfcPblancList.do?key=1004&cate_cd=1&searchCnd=40900000000&pageUnit=90"; d = "김포시육아종합지원센터 통진 분소 공공 실내놀이터 설계 및 제작・설치 용역 제안서 평가위원(후보자) 모집 공고"; e = 45784; f = 45786.369710648149 },
  @{ r = 79; a = "A120"; b = "충청도_태안군";   c = "http://eminwon.taean.go.kr/emwp/jsp/ofr/OfrNotAncmtL.jsp?not_ancmt_se_code=01,02,03,04,05&list_gubun=A"; d = "용역 입찰 공고(이원면 행정복지센터 건립사업 실시설계 용역)(제안공모)"; e = 45784; f = 45786.369710648149 },
  @{ r = 80; a = "A126"; b = "전라도_전주시";   c = "https://eminwon.jeonju.go.kr/emwp/jsp/ofr/OfrNotAncmtLSub.jsp?not_ancmt_se_code=01,02,03,04&epcCheck=Y&recent_mm=60&list_gubun=A"; d = "2025년 출연기관 경영평가 용역기관 선정 제안서 평가 결과 알림"; e = 45785; f = 45786.369710648149 },
  @{ r = 81; a = "A171"; b = "경상도_성주군";   c = "https://www.sj.go.kr/page.do?mnu_uid=1044&pageNo=1"; d = "2025년 수륜농협 참외AI 비파괴당도선별기 제작설치 제안서 평가위원(후보자) 모집 공고"; e = 45785; f = 45786.369710648149 },
  @{ r = 82; a = "A177"; b = "경상도_청도군";   c = "https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840"; d = "제안서 평가위원(후보자)모집 공고(화양읍 도시재생 뉴딜사업 지역역량강화용역)"; e = 45784; f = 45786.369710648149 },
  @{ r = 83; a = "A177"; b = "경상도_청도군";   c = "https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840"; d = "『온막천 소하천 정비사업 실시설계 용역』 신기술·특허공법 선정 기술제안서 제출"; e = 45784; f = 45786.369710648149 },
  @{ r = 84; a = "A177"; b = "경상도_청도군";   c = "https://www.cheongdo.go.kr/portal/saeol/gosi/list.do?mid=0301020000&token=1703813401840"; d = "제안서 평가위원(후보자)모집 공고(『온막천 소하천 정비사업 실시설계 용역』 특허공법 선정)"; e = 45784; f = 45786.369710648149 }
)

foreach ($row in $newRows) {
  $r = $row.r

  # Copy the date/time number-format (style index) from the row above so the
  # new E/F cells render as "yyyy-mm-dd hh:mm:ss" just like every other row,
  # instead of getting a brand-new style entry.
  $ws.Range("E$($r-1):F$($r-1)").Copy()
  $ws.Range("E${r}:F${r}").PasteSpecial(-4122)

  $ws.Cells.Item($r, 1).Value = $row.a
  $ws.Cells.Item($r, 2).Value = $row.b
  $ws.Cells.Item($r, 3).Value = $row.c
  $ws.Cells.Item($r, 4).Value = $row.d
  $ws.Cells.Item($r, 5).Value = $row.e
  $ws.Cells.Item($r, 6).Value = $row.f
}

# Match the author's final on-screen view: scrolled down near the new rows
# with I79 as the active cell.
$ws.Range("I79").Select()
$excel.ActiveWindow.ScrollRow = 70
